# Commit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# (Delete previous account-statement periods and add new ones; database updated)
#
# The underlying "Periodo Mora" list for the worker is reordered: what used to
# read ascending 2410,2411,2412,2501,2502,2503 (with 2504 fixed at the end)
# now reads descending 2503,2502,2501,2412,2411,2410 (2504 still fixed at the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2503"
$ws.Range("E17").Value = "2502"
$ws.Range("E18").Value = "2501"
$ws.Range("E19").Value = "2412"
$ws.Range("E20").Value = "2411"
$ws.Range("E21").Value = "2410"
$ws.Range("E22").Value = "2504"

$wb.Save()
